$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: TEXT "I knew it was at huzz" -> GOTO "JEFF_HUZZ" (condition unchanged) ---
$ws.Range("D11").Value = "JEFF_HUZZ"
$ws.Range("C11").Value = "GOTO"

# --- Row 12: END -> TEXT "IT WAS AT PLUZZ?" condition 2 ---
$ws.Range("C12").Value = "TEXT"
$ws.Range("D12").Value = "IT WAS AT PLUZZ?"
$ws.Range("E12").Value = 2

# --- Row 13: TEXT "IT WAS AT PLUZZ?" condition 2 -> END ---
$ws.Range("C13").Value = "END"
$ws.Range("D13").ClearContents()
$ws.Range("E13").ClearContents()

# --- Row 14: END -> TEXT "i didn't think it was at gruzz" condition 3 ---
$ws.Range("C14").Value = "TEXT"
$ws.Range("D14").Value = "i didn't think it was at gruzz"
$ws.Range("E14").Value = 3

# --- Row 15: TEXT "i didn't think it was at gruzz" condition 3 -> END ---
$ws.Range("C15").Value = "END"
$ws.Range("D15").ClearContents()
$ws.Range("E15").ClearContents()

# --- Row 16: END -> JEFF_HUZZ / TEXT (content filled in later) ---
$ws.Range("B16").Value = "JEFF_HUZZ"
$ws.Range("C16").Value = "TEXT"

# --- New row 17: CHOICE "But what's your name?[GRINGO,MORTIMER,GIBBY]" ---
$ws.Range("A17").Value = 15
$ws.Range("A17").HorizontalAlignment = -4131
$ws.Range("C17").Value = "CHOICE"
$ws.Range("D17").Value = "But what's your name?[GRINGO,MORTIMER,GIBBY]"

# --- New row 18: TEXT "your name is gringo?" condition 1 ---
$ws.Range("A18").Value = 16
$ws.Range("A18").HorizontalAlignment = -4131
$ws.Range("C18").Value = "TEXT"
$ws.Range("D18").Value = "your name is gringo?"
$ws.Range("E18").Value = 1

# --- New row 19: END ---
$ws.Range("A19").Value = 17
$ws.Range("A19").HorizontalAlignment = -4131
$ws.Range("C19").Value = "END"

# --- New row 20: TEXT "I knew your mortimer?" condition 2 ---
$ws.Range("A20").Value = 18
$ws.Range("A20").HorizontalAlignment = -4131
$ws.Range("C20").Value = "TEXT"
$ws.Range("D20").Value = "I knew your mortimer?"
$ws.Range("E20").Value = 2

# --- New row 21: END ---
$ws.Range("A21").Value = 19
$ws.Range("A21").HorizontalAlignment = -4131
$ws.Range("C21").Value = "END"

# --- New row 22: TEXT "bitch ass gibby" condition 3 ---
$ws.Range("A22").Value = 20
$ws.Range("A22").HorizontalAlignment = -4131
$ws.Range("C22").Value = "TEXT"
$ws.Range("D22").Value = "bitch ass gibby"
$ws.Range("E22").Value = 3

# --- Row 16 content filled in last ---
$ws.Range("D16").Value = "I forgot your name"

# --- View: zoom + selection on newly active cell ---
$excel.ActiveWindow.Zoom = 145
$ws.Range("D16").Select() | Out-Null
